$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Edits to the last three rows (tab-separated summary rows) ---
# Row 44: collapse the tab-delimited run list down to a single value.
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"

# Row 45: collapse the tab-delimited run list down to a single value.
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.01"

# Row 46: collapse the tab-delimited run list down to a single value.
$t.Rows.Item(46).Cells.Item(1).Range.Text = "175"

# --- Edits near the top of the table ---
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text = "202"
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00007"

# Row 9 (text "0.00002") is removed entirely.
$t.Rows.Item(9).Delete()

# After the deletion above, the former rows 11 and 12 are now rows 10 and 11.
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00004"
$t.Rows.Item(11).Cells.Item(1).Range.Text = "0.00004"

# Insert a brand-new row (text "0.00608") right after that row (now row 11),
# i.e. immediately before the row currently at index 12.
$newRow = $t.Rows.Add($t.Rows.Item(12))
$newRow.Cells.Item(1).Range.Text = "0.00608"
